$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cronjob re-sampled the "atypical sales" feed: the oldest day (2025-05-19,
# rows 2-6) dropped out of the window entirely, so remove those rows - the
# remaining rows shift up and the table shrinks from A1:H19 to A1:H14.
$ws.Range("A2:H6").EntireRow.Delete() | Out-Null

# Refresh the random row-order id in column A for every remaining data row.
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 4
$ws.Range("A4").Value = 6
$ws.Range("A5").Value = 9
$ws.Range("A6").Value = 10
$ws.Range("A7").Value = 11
$ws.Range("A8").Value = 1
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 7
$ws.Range("A11").Value = 2
$ws.Range("A12").Value = 5
$ws.Range("A13").Value = 8
$ws.Range("A14").Value = 12

# Within the unchanged "2025-05-26" group (rows 8-10) the three records were
# re-emitted in a different order; rotate B:H right by one row using a
# scratch area far outside the used range so cell types (shared-string /
# number / bool) are preserved exactly, then delete the scratch copy.
$ws.Range("B8:H10").Copy($ws.Range("B100:H102")) | Out-Null
$ws.Range("B100:H100").Copy($ws.Range("B9:H9")) | Out-Null
$ws.Range("B101:H101").Copy($ws.Range("B10:H10")) | Out-Null
$ws.Range("B102:H102").Copy($ws.Range("B8:H8")) | Out-Null
$ws.Range("B100:H102").Clear() | Out-Null

# Within the unchanged "2025-05-27" group, the last two records swapped
# places; same scratch-copy technique.
$ws.Range("B13:H13").Copy($ws.Range("B100:H100")) | Out-Null
$ws.Range("B14:H14").Copy($ws.Range("B13:H13")) | Out-Null
$ws.Range("B100:H100").Copy($ws.Range("B14:H14")) | Out-Null
$ws.Range("B100:H100").Clear() | Out-Null
